$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Group" values that changed for several students.
$ws.Range("C17").Value = 2
$ws.Range("C28").Value = 12
$ws.Range("C38").Value = 5
$ws.Range("C45").Value = 9
$ws.Range("C47").Value = 14
$ws.Range("C54").Value = 7

# Widen column A to fit the StudentID values (stores as width="22" in the xlsx).
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668

# Select row 21 (last interaction before save).
$ws.Range("A21:XFD21").Select()
